$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "'21.806.48"
$ws.Range("E2").Value = "  -1.61%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "'1.541.00"
$ws.Range("E3").Value = "  -1.24%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.55%  "

# Row 5 (USDC)
$ws.Range("E5").Value = "  +0.54%  "

# Row 6 (BNB)
$ws.Range("D6").Value = "'288.40"
$ws.Range("E6").Value = "  -0.14%  "

# Row 7 (XRP)
$ws.Range("D7").Value = "'0.3910"
$ws.Range("E7").Value = "  +2.88%  "

# Row 8 (Cardano)
$ws.Range("D8").Value = "'0.3200"
$ws.Range("E8").Value = "  -3.07%  "

# Row 9 (OKB)
$ws.Range("D9").Value = "'43.00"
$ws.Range("E9").Value = "  -1.84%  "

# Row 10 (Dogecoin)
$ws.Range("D10").Value = "'0.07148"
$ws.Range("E10").Value = "  -3.31%  "

# Row 11 (Polygon)
$ws.Range("D11").Value = "'1.059"
$ws.Range("E11").Value = "  -7.79%  "

# Row 12 (BinanceUSD)
$ws.Range("E12").Value = "  +0.57%  "

# Row 13 (Polkadot)
$ws.Range("E13").Value = "  -3.79%  "

# Row 14 (Solana)
$ws.Range("D14").Value = "'18.53"
$ws.Range("E14").Value = "  -8.40%  "

# Row 15 (Chainlink)
$ws.Range("D15").Value = "'6.617"
$ws.Range("E15").Value = "  -3.71%  "

# Row 16 (WrappedEther)
$ws.Range("D16").Value = "'1.550.49"
$ws.Range("E16").Value = "  -0.71%  "

# Row 17 (ShibaInu)
$ws.Range("E17").Value = "  -0.88%  "

# Row 18 (TRON)
$ws.Range("D18").Value = "'0.06562"
$ws.Range("E18").Value = "  -1.00%  "

# Row 19 (Litecoin)
$ws.Range("D19").Value = "'83.26"
$ws.Range("E19").Value = "  -3.42%  "

# Row 20 (Dai)
$ws.Range("E20").Value = "  +0.56%  "

# Row 21 (Uniswap)
$ws.Range("D21").Value = "'6.131"
$ws.Range("E21").Value = "  -4.49%  "

# Row 22 (Avalanche)
$ws.Range("D22").Value = "'15.23"
$ws.Range("E22").Value = "  -5.69%  "

# Row 23 (Cosmos)
$ws.Range("E23").Value = "  -6.20%  "

# Row 24 (Toncoin)
$ws.Range("D24").Value = "'2.399"
$ws.Range("E24").Value = "  +3.88%  "

# Row 25 (WrappedBTC)
$ws.Range("D25").Value = "'21.821.50"

# Row 26 (LidoDAOToken)
$ws.Range("D26").Value = "'2.378"
$ws.Range("E26").Value = "  -6.22%  "

# Row 27 (Monero)
$ws.Range("D27").Value = "'144.64"
$ws.Range("E27").Value = "  -3.64%  "

# Row 28 (EthereumClassic)
$ws.Range("D28").Value = "'18.35"
$ws.Range("E28").Value = "  -4.38%  "

# Row 29 (HuobiToken)
$ws.Range("D29").Value = "'4.858"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30 (WrappedliquidstakedEther2.0)
$ws.Range("D30").Value = "'1.722.99"
$ws.Range("E30").Value = "  -0.75%  "

# Row 31 (BitcoinCash)
$ws.Range("D31").Value = "'117.18"
$ws.Range("E31").Value = "  -3.95%  "

# Row 32 (ImmutableX)
$ws.Range("D32").Value = "'0.9630"
$ws.Range("E32").Value = "  -10.92%  "

# Row 33 (Filecoin)
$ws.Range("D33").Value = "'5.832"
$ws.Range("E33").Value = "  -2.19%  "

# Row 34 (Stellar)
$ws.Range("D34").Value = "'0.08212"
$ws.Range("E34").Value = "  -0.64%  "

# Row 35 (FraxShare)
$ws.Range("D35").Value = "'9.014"
$ws.Range("E35").Value = "  -3.63%  "

# Row 36 (WEMIXTOKEN)
$ws.Range("D36").Value = "'1.510"
$ws.Range("E36").Value = "  -18.40%  "

# Row 37 (Hedera)
$ws.Range("D37").Value = "'0.06087"
$ws.Range("E37").Value = "  -2.84%  "

# Row 38: was VeChain -> now InternetComputer(DFINITY)
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.088"
$ws.Range("E38").Value = "  -4.52%  "

# Row 39: was InternetComputer(DFINITY) -> now VeChain
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02222"
$ws.Range("E39").Value = "  -5.52%  "

# Row 40: was Algorand -> now TrustWalletToken
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'1.192"
$ws.Range("E40").Value = "  -4.99%  "

# Row 41: was TrustWalletToken -> now Algorand
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.2032"
$ws.Range("E41").Value = "  -6.51%  "

# Row 43 (Aptos)
$ws.Range("D43").Value = "'10.59"
$ws.Range("E43").Value = "  -4.64%  "

# Row 44 (TheSandbox)
$ws.Range("D44").Value = "'0.5719"
$ws.Range("E44").Value = "  -6.08%  "

# Row 45 (PancakeSwap)
$ws.Range("D45").Value = "'3.742"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46 (EnergySwap)
$ws.Range("D46").Value = "'12.82"
$ws.Range("E46").Value = "  -7.29%  "

# Row 47 (Decentraland)
$ws.Range("D47").Value = "'0.5513"
$ws.Range("E47").Value = "  -6.63%  "

# Row 48 (Quant)
$ws.Range("D48").Value = "'116.33"
$ws.Range("E48").Value = "  -4.97%  "

# Row 49 (NEARProtocol)
$ws.Range("D49").Value = "'1.851"
$ws.Range("E49").Value = "  -7.33%  "

# Row 50 (EOS)
$ws.Range("D50").Value = "'1.125"
$ws.Range("E50").Value = "  -4.57%  "

# Row 51 (Cronos)
$ws.Range("E51").Value = "  -4.05%  "
